$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "to meet; to see (a person) (person に)"
$ws.Range("A3").Value = "there is... (place に thing が)"
$ws.Range("A4").Value = "to buy (～を)"
$ws.Range("A5").Value = "to write (person に thing を)"
$ws.Range("A6").Value = "to take (a picture) (～を)"
$ws.Range("A7").Value = "to wait (～を)"
$ws.Range("A8").Value = "to understand (～が)"
$ws.Range("A9").Value = "(a person) is in...; stays at... (place に person が)"

$ws.Range("A40").Value = "right (～の)"
$ws.Range("A41").Value = "left (～の)"
$ws.Range("A42").Value = "front (～の)"
$ws.Range("A43").Value = "back (～の)"
$ws.Range("A44").Value = "inside (～の)"
$ws.Range("A45").Value = "on (～の)"
$ws.Range("A46").Value = "under (～の)"
$ws.Range("A47").Value = "near; nearby (～の)"
$ws.Range("A48").Value = "next (～の)"
$ws.Range("A49").Value = "between (A と B の)"
